$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit: Mars orbit altitude (D11) 384400 -> 1000 ---
$ws.Range("D11").Value = 1000

# --- New "Parent body" reference row (row 23) ---
$ws.Range("A23").Value = "Parent body"
$ws.Range("B23").Value = "-"
$ws.Range("C23").Value = "Earth"
$ws.Range("D23").Value = "Moon"
$ws.Range("E23").Value = "Mars"
$ws.Range("F23").Value = "Venus"
$ws.Range("G23").Value = "Europa"
$ws.Range("H23").Value = "Earth"

# Match the formatting used by the rest of the parameter table:
# columns A/B -> label font (Arial 10pt), C:F -> highlighted fill,
# G -> highlighted fill + 2-decimal number format, H -> plain bordered cell.
$ws.Range("A23:B23").Font.Name = "Arial"
$ws.Range("A23:B23").Font.Size = 10

$ws.Range("C23:F23").Interior.Pattern = 1
$ws.Range("C23:F23").Interior.Color = 11513855

$ws.Range("G23").Interior.Pattern = 1
$ws.Range("G23").Interior.Color = 11513855
$ws.Range("G23").NumberFormat = "0.00"

# --- Selection / view state ---
$ws.Range("D23").Select()
